# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500 ...
# Re-point the Alvearie/IBM StructureDefinition export at its new
# LinuxForHealth home, bump the version/date/publisher metadata, and
# drop the stale ele-1/ref-1 constraint text that no longer applies to
# the base Reference row.

$wb = $excel.ActiveWorkbook

# --- "Metadata" worksheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-with-period"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" worksheet -------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# referencePeriod extension (row 5, "Type(s)" column) now points at the
# linuxforhealth.org extension definition instead of ibm.com
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-period}`n"

# Reference (row 2, "Constraint(s)" column) no longer carries the
# ele-1/ref-1 constraint text
$elements.Range("AI2").Value = ""
